$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.9745285
$ws.Range("H2").Value = 13.949057
$ws.Range("I2").Value = 0.0400753517728026
$ws.Range("J2").Value = 0.02752929218001018
$ws.Range("M2").Value = 7.369448
$ws.Range("N2").Value = 14.738896
$ws.Range("O2").Value = 0.7452608427984224
$ws.Range("P2").Value = 0.661061693471796
$ws.Range("Q2").Value = 51.398425105268
$ws.Range("R2").Value = 205.593700421072
$ws.Range("S2").Value = 0.02986659043764212
$ws.Range("T2").Value = 0.0181985605085974
$ws.Range("G3").Value = 6.9745285
$ws.Range("H3").Value = 13.949057
$ws.Range("I3").Value = 0.0400753517728026
$ws.Range("J3").Value = 0.02752929218001018
$ws.Range("O3").Value = 0.01116592909756377
$ws.Range("P3").Value = 0.01485661309677453
$ws.Range("Q3").Value = 0.7700809401133333
$ws.Range("R3").Value = 4.620485640679999
$ws.Range("S3").Value = 0.0004474785364550405
$ws.Range("T3").Value = 0.0004089920427464718
$ws.Range("G4").Value = 6.9745285
$ws.Range("H4").Value = 13.949057
$ws.Range("I4").Value = 0.0400753517728026
$ws.Range("J4").Value = 0.02752929218001018
$ws.Range("M4").Value = 0.084843
$ws.Range("N4").Value = 0.254529
$ws.Range("O4").Value = 0.008580040959044227
$ws.Range("P4").Value = 0.0114160091622658
$ws.Range("Q4").Value = 0.5917399215255
$ws.Range("R4").Value = 3.550439529153
$ws.Range("S4").Value = 0.000343848159658752
$ws.Range("T4").Value = 0.0003142746517576885
$ws.Range("G5").Value = 6.9745285
$ws.Range("H5").Value = 13.949057
$ws.Range("I5").Value = 0.0400753517728026
$ws.Range("J5").Value = 0.02752929218001018
$ws.Range("M5").Value = 2.32371
$ws.Range("N5").Value = 6.97113
$ws.Range("O5").Value = 0.2349931871449696
$ws.Range("P5").Value = 0.3126656842691638
$ws.Range("Q5").Value = 16.206781620735
$ws.Range("R5").Value = 97.24068972441
$ws.Range("S5").Value = 0.009417434639046693
$ws.Range("T5").Value = 0.008607464976908624
$ws.Range("I6").Value = 0.2331362278651957
$ws.Range("J6").Value = 0.2402252900123789
$ws.Range("M6").Value = 7.369448
$ws.Range("N6").Value = 14.738896
$ws.Range("O6").Value = 0.7452608427984224
$ws.Range("P6").Value = 0.661061693471796
$ws.Range("Q6").Value = 299.0076048536694
$ws.Range("R6").Value = 1794.045629122016
$ws.Range("S6").Value = 0.1737473016656608
$ws.Range("T6").Value = 0.1588037370303365
$ws.Range("I7").Value = 0.2331362278651957
$ws.Range("J7").Value = 0.2402252900123789
$ws.Range("O7").Value = 0.01116592909756377
$ws.Range("P7").Value = 0.01485661309677453
$ws.Range("S7").Value = 0.002603182590416247
$ws.Range("T7").Value = 0.003568934189774367
$ws.Range("I8").Value = 0.2331362278651957
$ws.Range("J8").Value = 0.2402252900123789
$ws.Range("M8").Value = 0.084843
$ws.Range("N8").Value = 0.254529
$ws.Range("O8").Value = 0.008580040959044227
$ws.Range("P8").Value = 0.0114160091622658
$ws.Range("Q8").Value = 3.442415526726
$ws.Range("R8").Value = 30.981739740534
$ws.Range("S8").Value = 0.002000318384120447
$ws.Range("T8").Value = 0.002742414111789276
$ws.Range("I9").Value = 0.2331362278651957
$ws.Range("J9").Value = 0.2402252900123789
$ws.Range("M9").Value = 2.32371
$ws.Range("N9").Value = 6.97113
$ws.Range("O9").Value = 0.2349931871449696
$ws.Range("P9").Value = 0.3126656842691638
$ws.Range("Q9").Value = 94.28209025622
$ws.Range("R9").Value = 848.5388123059801
$ws.Range("S9").Value = 0.05478542522499823
$ws.Range("T9").Value = 0.07511020468047877
$ws.Range("G10").Value = 16.42157633333333
$ws.Range("H10").Value = 49.264729
$ws.Range("I10").Value = 0.09435769718659309
$ws.Range("J10").Value = 0.09722686765205855
$ws.Range("M10").Value = 7.369448
$ws.Range("N10").Value = 14.738896
$ws.Range("O10").Value = 0.7452608427984224
$ws.Range("P10").Value = 0.661061693471796
$ws.Range("Q10").Value = 121.0179528665307
$ws.Range("R10").Value = 726.1077171991841
$ws.Range("S10").Value = 0.0703210969297987
$ws.Range("T10").Value = 0.06427295778102801
$ws.Range("G11").Value = 16.42157633333333
$ws.Range("H11").Value = 49.264729
$ws.Range("I11").Value = 0.09435769718659309
$ws.Range("J11").Value = 0.09722686765205855
$ws.Range("O11").Value = 0.01116592909756377
$ws.Range("P11").Value = 0.01485661309677453
$ws.Range("Q11").Value = 1.813160981551111
$ws.Range("R11").Value = 16.31844883396
$ws.Range("S11").Value = 0.001053591356594891
$ws.Range("T11").Value = 0.001444461955317937
$ws.Range("G12").Value = 16.42157633333333
$ws.Range("H12").Value = 49.264729
$ws.Range("I12").Value = 0.09435769718659309
$ws.Range("J12").Value = 0.09722686765205855
$ws.Range("M12").Value = 0.084843
$ws.Range("N12").Value = 0.254529
$ws.Range("O12").Value = 0.008580040959044227
$ws.Range("P12").Value = 0.0114160091622658
$ws.Range("Q12").Value = 1.393255800849
$ws.Range("R12").Value = 12.539302207641
$ws.Range("S12").Value = 0.0008095929066620609
$ws.Range("T12").Value = 0.001109942811934305
$ws.Range("G13").Value = 16.42157633333333
$ws.Range("H13").Value = 49.264729
$ws.Range("I13").Value = 0.09435769718659309
$ws.Range("J13").Value = 0.09722686765205855
$ws.Range("M13").Value = 2.32371
$ws.Range("N13").Value = 6.97113
$ws.Range("O13").Value = 0.2349931871449696
$ws.Range("P13").Value = 0.3126656842691638
$ws.Range("Q13").Value = 38.15898114153001
$ws.Range("R13").Value = 343.43083027377
$ws.Range("S13").Value = 0.02217341599353745
$ws.Range("T13").Value = 0.03039950510377832
$ws.Range("G14").Value = 8.432852499999999
$ws.Range("H14").Value = 16.865705
$ws.Range("I14").Value = 0.04845482105143851
$ws.Range("J14").Value = 0.03328547017671937
$ws.Range("M14").Value = 7.369448
$ws.Range("N14").Value = 14.738896
$ws.Range("O14").Value = 0.7452608427984224
$ws.Range("P14").Value = 0.661061693471796
$ws.Range("Q14").Value = 62.14546799041999
$ws.Range("R14").Value = 248.58187196168
$ws.Range("S14").Value = 0.0361114807744418
$ws.Range("T14").Value = 0.02200374928302707
$ws.Range("G15").Value = 8.432852499999999
$ws.Range("H15").Value = 16.865705
$ws.Range("I15").Value = 0.04845482105143851
$ws.Range("J15").Value = 0.03328547017671937
$ws.Range("O15").Value = 0.01116592909756377
$ws.Range("P15").Value = 0.01485661309677453
$ws.Range("Q15").Value = 0.9310993540333331
$ws.Range("R15").Value = 5.586596124199999
$ws.Range("S15").Value = 0.0005410430962955028
$ws.Range("T15").Value = 0.000494509352159747
$ws.Range("G16").Value = 8.432852499999999
$ws.Range("H16").Value = 16.865705
$ws.Range("I16").Value = 0.04845482105143851
$ws.Range("J16").Value = 0.03328547017671937
$ws.Range("M16").Value = 0.084843
$ws.Range("N16").Value = 0.254529
$ws.Range("O16").Value = 0.008580040959044227
$ws.Range("P16").Value = 0.0114160091622658
$ws.Range("Q16").Value = 0.7154685046574999
$ws.Range("R16").Value = 4.292811027945
$ws.Range("S16").Value = 0.0004157443492845009
$ws.Range("T16").Value = 0.0003799872325077534
$ws.Range("G17").Value = 8.432852499999999
$ws.Range("H17").Value = 16.865705
$ws.Range("I17").Value = 0.04845482105143851
$ws.Range("J17").Value = 0.03328547017671937
$ws.Range("M17").Value = 2.32371
$ws.Range("N17").Value = 6.97113
$ws.Range("O17").Value = 0.2349931871449696
$ws.Range("P17").Value = 0.3126656842691638
$ws.Range("Q17").Value = 19.595503682775
$ws.Range("R17").Value = 117.57302209665
$ws.Range("S17").Value = 0.0113865528314167
$ws.Range("T17").Value = 0.01040722430902481
$ws.Range("G18").Value = 13.375494
$ws.Range("H18").Value = 40.126482
$ws.Range("I18").Value = 0.07685503431307371
$ws.Range("J18").Value = 0.07919199463690765
$ws.Range("M18").Value = 7.369448
$ws.Range("N18").Value = 14.738896
$ws.Range("O18").Value = 0.7452608427984224
$ws.Range("P18").Value = 0.661061693471796
$ws.Range("Q18").Value = 98.57000750731198
$ws.Range("R18").Value = 591.4200450438719
$ws.Range("S18").Value = 0.05727704764546299
$ws.Range("T18").Value = 0.05235079408408356
$ws.Range("G19").Value = 13.375494
$ws.Range("H19").Value = 40.126482
$ws.Range("I19").Value = 0.07685503431307371
$ws.Range("J19").Value = 0.07919199463690765
$ws.Range("O19").Value = 0.01116592909756377
$ws.Range("P19").Value = 0.01485661309677453
$ws.Range("Q19").Value = 1.47683287752
$ws.Range("R19").Value = 13.29149589768
$ws.Range("S19").Value = 0.0008581578639306117
$ws.Range("T19").Value = 0.00117652482468238
$ws.Range("G20").Value = 13.375494
$ws.Range("H20").Value = 40.126482
$ws.Range("I20").Value = 0.07685503431307371
$ws.Range("J20").Value = 0.07919199463690765
$ws.Range("M20").Value = 0.084843
$ws.Range("N20").Value = 0.254529
$ws.Range("O20").Value = 0.008580040959044227
$ws.Range("P20").Value = 0.0114160091622658
$ws.Range("Q20").Value = 1.134817037442
$ws.Range("R20").Value = 10.213353336978
$ws.Range("S20").Value = 0.0006594193423149219
$ws.Range("T20").Value = 0.0009040565363530419
$ws.Range("G21").Value = 13.375494
$ws.Range("H21").Value = 40.126482
$ws.Range("I21").Value = 0.07685503431307371
$ws.Range("J21").Value = 0.07919199463690765
$ws.Range("M21").Value = 2.32371
$ws.Range("N21").Value = 6.97113
$ws.Range("O21").Value = 0.2349931871449696
$ws.Range("P21").Value = 0.3126656842691638
$ws.Range("Q21").Value = 31.08076916274
$ws.Range("R21").Value = 279.72692246466
$ws.Range("S21").Value = 0.01806040946136519
$ws.Range("T21").Value = 0.02476061919178868
$ws.Range("G22").Value = 88.25696566666666
$ws.Range("H22").Value = 264.770897
$ws.Range("I22").Value = 0.5071208678108963
$ws.Range("J22").Value = 0.5225410853419253
$ws.Range("M22").Value = 7.369448
$ws.Range("N22").Value = 14.738896
$ws.Range("O22").Value = 0.7452608427984224
$ws.Range("P22").Value = 0.661061693471796
$ws.Range("Q22").Value = 650.4051191182853
$ws.Range("R22").Value = 3902.430714709712
$ws.Range("S22").Value = 0.3779373253454159
$ws.Range("T22").Value = 0.3454318947847234
$ws.Range("G23").Value = 88.25696566666666
$ws.Range("H23").Value = 264.770897
$ws.Range("I23").Value = 0.5071208678108963
$ws.Range("J23").Value = 0.5225410853419253
$ws.Range("O23").Value = 0.01116592909756377
$ws.Range("P23").Value = 0.01485661309677453
$ws.Range("Q23").Value = 9.74474576914222
$ws.Range("R23").Value = 87.70271192227999
$ws.Range("S23").Value = 0.005662475653871477
$ws.Range("T23").Value = 0.007763190732093624
$ws.Range("G24").Value = 88.25696566666666
$ws.Range("H24").Value = 264.770897
$ws.Range("I24").Value = 0.5071208678108963
$ws.Range("J24").Value = 0.5225410853419253
$ws.Range("M24").Value = 0.084843
$ws.Range("N24").Value = 0.254529
$ws.Range("O24").Value = 0.008580040959044227
$ws.Range("P24").Value = 0.0114160091622658
$ws.Range("Q24").Value = 7.487985738057
$ws.Range("R24").Value = 67.391871642513
$ws.Range("S24").Value = 0.004351117817003543
$ws.Range("T24").Value = 0.005965333817923735
$ws.Range("G25").Value = 88.25696566666666
$ws.Range("H25").Value = 264.770897
$ws.Range("I25").Value = 0.5071208678108963
$ws.Range("J25").Value = 0.5225410853419253
$ws.Range("M25").Value = 2.32371
$ws.Range("N25").Value = 6.97113
$ws.Range("O25").Value = 0.2349931871449696
$ws.Range("P25").Value = 0.3126656842691638
$ws.Range("Q25").Value = 205.08359368929
$ws.Range("R25").Value = 1845.75234320361
$ws.Range("S25").Value = 0.1191699489946054
$ws.Range("T25").Value = 0.1633806660071846
